$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("NextBus1")
$ws.Cells.Item(2, 6).Value = 45688.39193287037
$ws.Cells.Item(2, 12).Value = "SD"
$ws.Cells.Item(2, 15).Value = 7
$ws.Cells.Item(3, 6).Value = 45688.39123842592
$ws.Cells.Item(3, 15).Value = 6
$ws.Cells.Item(4, 6).Value = 45688.38809027777
$ws.Cells.Item(4, 15).Value = 1
$ws.Cells.Item(5, 6).Value = 45688.39040509259
$ws.Cells.Item(5, 15).Value = 4
$ws.Cells.Item(6, 6).Value = 45688.39842592592
$ws.Cells.Item(6, 15).Value = 16
$ws.Cells.Item(7, 6).Value = 45688.39006944445
$ws.Cells.Item(7, 15).Value = 4
$ws.Cells.Item(8, 6).Value = 45688.39489583333
$ws.Cells.Item(8, 15).Value = 11
$ws.Cells.Item(9, 6).Value = 45688.39636574074
$ws.Cells.Item(9, 9).Value = "SDA"
$ws.Cells.Item(9, 15).Value = 13
$ws.Cells.Item(10, 6).Value = 45688.39023148148
$ws.Cells.Item(10, 15).Value = 4
$ws.Cells.Item(11, 6).Value = 45688.39292824074
$ws.Cells.Item(11, 15).Value = 8
$ws.Cells.Item(12, 6).Value = 45688.38905092593
$ws.Cells.Item(12, 12).Value = "SD"
$ws.Cells.Item(12, 15).Value = 3
$ws.Cells.Item(13, 6).Value = 45688.39908564815
$ws.Cells.Item(13, 12).Value = "SD"
$ws.Cells.Item(13, 15).Value = 17
$ws.Cells.Item(14, 6).Value = 45688.39046296296
$ws.Cells.Item(14, 15).Value = 5
$ws.Cells.Item(15, 6).Value = 45688.39512731481
$ws.Cells.Item(15, 15).Value = 11

$ws = $wb.Worksheets.Item("NextBus2")
$ws.Cells.Item(2, 6).Value = 45688.4024074074
$ws.Cells.Item(2, 15).Value = 22
$ws.Cells.Item(3, 6).Value = 45688.39899305555
$ws.Cells.Item(3, 15).Value = 17
$ws.Cells.Item(4, 6).Value = 45688.39193287037
$ws.Cells.Item(4, 15).Value = 7
$ws.Cells.Item(5, 6).Value = 45688.39825231482
$ws.Cells.Item(5, 15).Value = 16
$ws.Cells.Item(6, 6).Value = 45688.40359953704
$ws.Cells.Item(6, 15).Value = 23
$ws.Cells.Item(7, 6).Value = 45688.39987268519
$ws.Cells.Item(7, 15).Value = 18
$ws.Cells.Item(8, 6).Value = 45688.39704861111
$ws.Cells.Item(8, 15).Value = 14
$ws.Cells.Item(9, 6).Value = 45688.40331018518
$ws.Cells.Item(9, 9).Value = "SEA"
$ws.Cells.Item(9, 15).Value = 23
$ws.Cells.Item(10, 6).Value = 45688.39284722223
$ws.Cells.Item(10, 15).Value = 8
$ws.Cells.Item(11, 15).Value = 23
$ws.Cells.Item(12, 6).Value = 45688.39483796297
$ws.Cells.Item(12, 15).Value = 11
$ws.Cells.Item(13, 6).Value = 45688.40502314815
$ws.Cells.Item(13, 12).Value = "DD"
$ws.Cells.Item(13, 15).Value = 26
$ws.Cells.Item(14, 6).Value = 45688.39288194444
$ws.Cells.Item(14, 15).Value = 8
$ws.Cells.Item(15, 6).Value = 45688.39634259259
$ws.Cells.Item(15, 15).Value = 13

$ws = $wb.Worksheets.Item("NextBus3")
$ws.Cells.Item(2, 6).Value = 45688.40834490741
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 15).Value = 30
$ws.Cells.Item(3, 6).Value = 45688.40403935185
$ws.Cells.Item(3, 15).Value = 24
$ws.Cells.Item(4, 6).Value = 45688.39429398148
$ws.Cells.Item(4, 12).Value = "DD"
$ws.Cells.Item(4, 15).Value = 10
$ws.Cells.Item(5, 6).Value = 45688.40501157408
$ws.Cells.Item(5, 15).Value = 26
$ws.Cells.Item(6, 6).Value = 45688.41496527778
$ws.Cells.Item(6, 15).Value = 40
$ws.Cells.Item(7, 6).Value = 45688.40578703704
$ws.Cells.Item(7, 12).Value = "DD"
$ws.Cells.Item(7, 15).Value = 27
$ws.Cells.Item(8, 6).Value = 45688.40618055555
$ws.Cells.Item(8, 15).Value = 27
$ws.Cells.Item(9, 6).Value = 45688.41425925926
$ws.Cells.Item(9, 15).Value = 39
$ws.Cells.Item(10, 6).Value = 45688.40679398148
$ws.Cells.Item(10, 15).Value = 28
$ws.Cells.Item(11, 15).Value = 38
$ws.Cells.Item(12, 6).Value = 45688.40255787037
$ws.Cells.Item(12, 12).Value = "DD"
$ws.Cells.Item(12, 15).Value = 22
$ws.Cells.Item(13, 6).Value = 45688.41442129629
$ws.Cells.Item(13, 12).Value = "SD"
$ws.Cells.Item(13, 15).Value = 39
$ws.Cells.Item(14, 15).Value = 20
$ws.Cells.Item(15, 6).Value = 45688.40453703704
$ws.Cells.Item(15, 15).Value = 25
